# "update da lista de peças"
# Renumbers/edits the screws list at the bottom of the parts table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: note changes from "3 peças, falta botão stopper" to "falta definir comprimento"
$ws.Range("G31").Value = "falta definir comprimento"

# Row 32: "parafuso CHC M4x38" x24 -> "parafuso CHC M4x40" x28
$ws.Range("A32").Value = "parafuso CHC M4x40"
$ws.Range("B32").Value = 28

# Row 33: "parafuso CHC M3x" -> "parafuso CHC M3x12", note -> "3 peças, falta botão stopper"
$ws.Range("A33").Value = "parafuso CHC M3x12"
$ws.Range("G33").Value = "3 peças, falta botão stopper"

# Row 34: "parafuso CHC M5x47" x4 -> "parafuso CHC M5x50" x4
$ws.Range("A34").Value = "parafuso CHC M5x50"

# Row 35 stays "parafuso CHC M3x" x3 - values unchanged

# New row 36: parafuso CHC M5x35
$ws.Range("A36").Value = "parafuso CHC M5x35"
$ws.Range("B36").Value = 2
$ws.Range("C36").Value = "comprar"
$ws.Range("E36").Value = "não"

# Row 37 content moves from the blank trailing row to a real data row,
# keeping F37's existing (empty) styled cell.
$ws.Range("A37").Value = "parafuso CHC M4x20"
$ws.Range("B37").Value = 2
$ws.Range("C37").Value = "comprar"
$ws.Range("E37").Value = "não"

# New row 38: parafuso M3xl (l<10)
$ws.Range("A38").Value = "parafuso M3xl (l<10)"
$ws.Range("B38").Value = 4
$ws.Range("C38").Value = "comprar"
$ws.Range("E38").Value = "não"

# Update sheet view: selection moves to E35:E38, no frozen/scrolled top-left cell.
$ws.Range("E35:E38").Select()
